$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HMI Internal")

$rows = @(
    ,@(12, "SEC0_LOCK_SEMI_PART", "BIT", 1)
    ,@(13, "SEC0_LOCK_MANU_PART", "BIT", 1)
    ,@(14, "SEC1_LOCK_SEMI_PART", "BIT", 1)
    ,@(15, "SEC1_LOCK_MANU_PART", "BIT", 1)
    ,@(16, "SEC2_LOCK_SEMI_PART", "BIT", 1)
    ,@(17, "SEC2_LOCK_MANU_PART", "BIT", 1)
    ,@(18, "SHELF_AUTO_WATER_REG_SIZE", "WORD", 8)
    ,@(19, "SHELF_AUTO_CURRENT_NO", "WORD", 1)
    ,@(20, "SHELF_AUTO_PLC_OFFSET_ADDR", "WORD", 1)
    ,@(21, "SHELF_AUTO_OFFSET_ADDR", "WORD", 1)
    ,@(22, "S0_WATER_ON_HH", "WORD", 1)
    ,@(23, "S0_WATER_ON_MM", "WORD", 1)
    ,@(24, "S0_WATER_OFF_HH", "WORD", 1)
    ,@(25, "S0_WATER_OFF_MM", "WORD", 1)
    ,@(26, "S0_LIGHT_ON_HH", "WORD", 1)
    ,@(27, "S0_LIGHT_ON_MM", "WORD", 1)
    ,@(28, "S0_LIGHT_OFF_HH", "WORD", 1)
    ,@(29, "S0_LIGHT_OFF_MM", "WORD", 1)
    ,@(30, "S1_WATER_ON_HH", "WORD", 1)
    ,@(31, "S1_WATER_ON_MM", "WORD", 1)
    ,@(32, "S1_WATER_OFF_HH", "WORD", 1)
    ,@(33, "S1_WATER_OFF_MM", "WORD", 1)
    ,@(34, "S1_LIGHT_ON_HH", "WORD", 1)
    ,@(35, "S1_LIGHT_ON_MM", "WORD", 1)
    ,@(36, "S1_LIGHT_OFF_HH", "WORD", 1)
    ,@(37, "S1_LIGHT_OFF_MM", "WORD", 1)
    ,@(38, "S2_WATER_ON_HH", "WORD", 1)
    ,@(39, "S2_WATER_ON_MM", "WORD", 1)
    ,@(40, "S2_WATER_OFF_HH", "WORD", 1)
    ,@(41, "S2_WATER_OFF_MM", "WORD", 1)
    ,@(42, "S2_LIGHT_ON_HH", "WORD", 1)
    ,@(43, "S2_LIGHT_ON_MM", "WORD", 1)
    ,@(44, "S2_LIGHT_OFF_HH", "WORD", 1)
    ,@(45, "S2_LIGHT_OFF_MM", "WORD", 1)
    ,@(46, "S3_WATER_ON_HH", "WORD", 1)
    ,@(47, "S3_WATER_ON_MM", "WORD", 1)
    ,@(48, "S3_WATER_OFF_HH", "WORD", 1)
    ,@(49, "S3_WATER_OFF_MM", "WORD", 1)
    ,@(50, "S3_LIGHT_ON_HH", "WORD", 1)
    ,@(51, "S3_LIGHT_ON_MM", "WORD", 1)
    ,@(52, "S3_LIGHT_OFF_HH", "WORD", 1)
    ,@(53, "S3_LIGHT_OFF_MM", "WORD", 1)
    ,@(54, "S4_WATER_ON_HH", "WORD", 1)
    ,@(55, "S4_WATER_ON_MM", "WORD", 1)
    ,@(56, "S4_WATER_OFF_HH", "WORD", 1)
    ,@(57, "S4_WATER_OFF_MM", "WORD", 1)
    ,@(58, "S4_LIGHT_ON_HH", "WORD", 1)
    ,@(59, "S4_LIGHT_ON_MM", "WORD", 1)
    ,@(60, "S4_LIGHT_OFF_HH", "WORD", 1)
    ,@(61, "S4_LIGHT_OFF_MM", "WORD", 1)
    ,@(62, "S5_WATER_ON_HH", "WORD", 1)
    ,@(63, "S5_WATER_ON_MM", "WORD", 1)
    ,@(64, "S5_WATER_OFF_HH", "WORD", 1)
    ,@(65, "S5_WATER_OFF_MM", "WORD", 1)
    ,@(66, "S5_LIGHT_ON_HH", "WORD", 1)
    ,@(67, "S5_LIGHT_ON_MM", "WORD", 1)
    ,@(68, "S5_LIGHT_OFF_HH", "WORD", 1)
    ,@(69, "S5_LIGHT_OFF_MM", "WORD", 1)
    ,@(70, "S6_WATER_ON_HH", "WORD", 1)
    ,@(71, "S6_WATER_ON_MM", "WORD", 1)
    ,@(72, "S6_WATER_OFF_HH", "WORD", 1)
    ,@(73, "S6_WATER_OFF_MM", "WORD", 1)
    ,@(74, "S6_LIGHT_ON_HH", "WORD", 1)
    ,@(75, "S6_LIGHT_ON_MM", "WORD", 1)
    ,@(76, "S6_LIGHT_OFF_HH", "WORD", 1)
    ,@(77, "S6_LIGHT_OFF_MM", "WORD", 1)
    ,@(78, "S7_WATER_ON_HH", "WORD", 1)
    ,@(79, "S7_WATER_ON_MM", "WORD", 1)
    ,@(80, "S7_WATER_OFF_HH", "WORD", 1)
    ,@(81, "S7_WATER_OFF_MM", "WORD", 1)
    ,@(82, "S7_LIGHT_ON_HH", "WORD", 1)
    ,@(83, "S7_LIGHT_ON_MM", "WORD", 1)
    ,@(84, "S7_LIGHT_OFF_HH", "WORD", 1)
    ,@(85, "S7_LIGHT_OFF_MM", "WORD", 1)
    ,@(86, "S8_WATER_ON_HH", "WORD", 1)
    ,@(87, "S8_WATER_ON_MM", "WORD", 1)
    ,@(88, "S8_WATER_OFF_HH", "WORD", 1)
    ,@(89, "S8_WATER_OFF_MM", "WORD", 1)
    ,@(90, "S8_LIGHT_ON_HH", "WORD", 1)
    ,@(91, "S8_LIGHT_ON_MM", "WORD", 1)
    ,@(92, "S8_LIGHT_OFF_HH", "WORD", 1)
    ,@(93, "S8_LIGHT_OFF_MM", "WORD", 1)
    ,@(94, "S9_WATER_ON_HH", "WORD", 1)
    ,@(95, "S9_WATER_ON_MM", "WORD", 1)
    ,@(96, "S9_WATER_OFF_HH", "WORD", 1)
    ,@(97, "S9_WATER_OFF_MM", "WORD", 1)
    ,@(98, "S9_LIGHT_ON_HH", "WORD", 1)
    ,@(99, "S9_LIGHT_ON_MM", "WORD", 1)
    ,@(100, "S9_LIGHT_OFF_HH", "WORD", 1)
    ,@(101, "S9_LIGHT_OFF_MM", "WORD", 1)
    ,@(102, "S0_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(103, "S1_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(104, "S2_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(105, "S3_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(106, "S4_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(107, "S5_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(108, "S6_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(109, "S7_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(110, "S8_PLC_ADDR_OFFSET", "WORD", 1)
    ,@(111, "S9_PLC_ADDR_OFFSET", "WORD", 1)
)

foreach ($row in $rows) {
    $r = $row[0]
    $name = $row[1]
    $typ = $row[2]
    $val = $row[3]
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $typ
    $ws.Cells.Item($r, 4).Value = $val
}

$ws.Activate()
$ws.Range("F13").Select()
